$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Police Commissioner name, volume/issue number, report week dates) ---
$ws.Range("M6").Value = "Thomas G. Donlon"
$ws.Range("A8").Characters(21, 2).Text = "39"
$ws.Range("C9").Characters(27, 9).Text = "9/23/2024"
$ws.Range("C9").Characters(47, 9).Text = "9/29/2024"

# --- Weekly crime-complaint statistics table (rows 15-31) ---
# Row 15
$ws.Range("C15").Value = 1
$ws.Range("I16").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("F15").Value = 1
$ws.Range("I16").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = -14.285714285714
$ws.Range("L15").Value = -57.142857142857
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -72.727272727272

# Row 16
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -23.076923076923
$ws.Range("I16").Value = 96
$ws.Range("J16").Value = 110
$ws.Range("K16").Value = -12.727272727272
$ws.Range("L16").Value = 18.518518518518
$ws.Range("M16").Value = -21.951219512195
$ws.Range("N16").Value = -83.476764199655

# Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = -28.571428571428
$ws.Range("I17").Value = 134
$ws.Range("J17").Value = 137
$ws.Range("K17").Value = -2.189781021897
$ws.Range("L17").Value = 15.517241379310
$ws.Range("M17").Value = 48.888888888888
$ws.Range("N17").Value = -38.812785388127

# Row 18
$ws.Range("C18").Value = "'0"
$ws.Range("D15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = "'0"
$ws.Range("D15").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "'***.*"
$ws.Range("D15").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -77.777777777777
$ws.Range("L18").Value = -57.291666666666
$ws.Range("M18").Value = -30.508474576271
$ws.Range("N18").Value = -91.067538126361

# Row 19
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 25
$ws.Range("H19").Value = -24.242424242424
$ws.Range("I19").Value = 252
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = -16
$ws.Range("L19").Value = -1.945525291828
$ws.Range("M19").Value = 25.373134328358
$ws.Range("N19").Value = -49.900596421471

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = -57.142857142857
$ws.Range("I20").Value = 29
$ws.Range("J20").Value = 82
$ws.Range("K20").Value = -64.634146341463
$ws.Range("L20").Value = -47.272727272727
$ws.Range("M20").Value = 45
$ws.Range("N20").Value = -89.860139860139

# Row 21
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -23.529411764705
$ws.Range("F21").Value = 51
$ws.Range("G21").Value = 76
$ws.Range("H21").Value = -32.894736842105
$ws.Range("I21").Value = 559
$ws.Range("J21").Value = 709
$ws.Range("K21").Value = -21.156558533145
$ws.Range("L21").Value = -9.838709677419
$ws.Range("M21").Value = 10.256410256410
$ws.Range("N21").Value = -73.112073112073

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("I16").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").Value = "'0"
$ws.Range("D15").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("D15").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 2
$ws.Range("I16").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 16
$ws.Range("K22").Value = -40.740740740740
$ws.Range("L22").Value = 14.285714285714
$ws.Range("M22").Value = -30.434782608695

# Row 23
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -11.111111111111
$ws.Range("I23").Value = 99
$ws.Range("J23").Value = 114
$ws.Range("K23").Value = -13.157894736842
$ws.Range("L23").Value = 11.235955056179
$ws.Range("M23").Value = 59.677419354838

# Row 24
$ws.Range("C24").Value = 4
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = -50
$ws.Range("F24").Value = 35
$ws.Range("H24").Value = -28.571428571428
$ws.Range("I24").Value = 308
$ws.Range("J24").Value = 386
$ws.Range("K24").Value = -20.207253886010
$ws.Range("L24").Value = -20.207253886010
$ws.Range("M24").Value = -36.099585062240

# Row 25
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -87.5
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = -36.842105263157
$ws.Range("I25").Value = 52
$ws.Range("J25").Value = 132
$ws.Range("K25").Value = -60.606060606060
$ws.Range("L25").Value = -61.194029850746

# Row 26
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 15
$ws.Range("H26").Value = 73.333333333333
$ws.Range("I26").Value = 181
$ws.Range("J26").Value = 181
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 10.365853658536
$ws.Range("M26").Value = -19.196428571428

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("I16").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("I16").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = 0
$ws.Range("K16").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 1
$ws.Range("I16").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("G27").Value = 1
$ws.Range("I16").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("H27").Value = 0
$ws.Range("K16").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("I27").Value = 11
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = 22.222222222222
$ws.Range("L27").Value = -45

# Row 28
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -40
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -37.5
$ws.Range("I28").Value = 31
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = -13.888888888888
$ws.Range("L28").Value = -13.888888888888

# Row 29
$ws.Range("N29").Value = -76

# Row 30
$ws.Range("N30").Value = -76.190476190476

# Row 31
$ws.Range("C31").Value = 1
$ws.Range("I16").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 1
$ws.Range("I16").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("H31").Value = -50
$ws.Range("I31").Value = 8
$ws.Range("J31").Value = 5
$ws.Range("K31").Value = 60
$ws.Range("L31").Value = -11.111111111111
